$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to remain plain text,
# matching the source data which stores them as inline strings
# (many values, e.g. "243.80" or "0.9998", look numeric and would
# otherwise be auto-converted by Excel on assignment).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.390.13"

$ws.Range("D3").Value = "1.868.57"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "243.80"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").Value = "0.7044"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").Value = "0.3133"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").Value = "24.46"
$ws.Range("E10").Value = "  -2.14%  "

$ws.Range("D11").Value = "0.07848"
$ws.Range("E11").Value = "  -4.64%  "

$ws.Range("D12").Value = "1.862.12"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").Value = "93.74"
$ws.Range("E13").Value = "  -1.17%  "

$ws.Range("D14").Value = "5.180"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "0.7010"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").Value = "6.515"
$ws.Range("E16").Value = "  +1.59%  "

$ws.Range("D17").Value = "0.000008398"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "29.379.20"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "254.12"
$ws.Range("E19").Value = "  +4.25%  "

$ws.Range("D20").Value = "2.117.90"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("D21").Value = "13.10"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").Value = "0.9998"

$ws.Range("D23").Value = "7.645"
$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D26").Value = "9.006"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").Value = "161.06"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("D29").Value = "1.500"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").Value = "4.318"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").Value = "4.254"
$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").Value = "1.213"
$ws.Range("E32").Value = "  +3.18%  "

$ws.Range("D33").Value = "0.05281"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.173"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7483"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  +0.88%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01878"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.280.50"
$ws.Range("E39").Value = "  +1.22%  "

$ws.Range("D40").Value = "2.765"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").Value = "0.8923"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("D42").Value = "108.70"
$ws.Range("E42").Value = "  -3.57%  "

$ws.Range("D43").Value = "5.998"
$ws.Range("E43").Value = "  -7.24%  "

$ws.Range("D44").Value = "70.94"
$ws.Range("E44").Value = "  -4.22%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  -3.55%  "

$ws.Range("D47").Value = "2.017.97"
$ws.Range("E47").Value = "  -1.01%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.600"
$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.798"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").Value = "0.4301"
$ws.Range("E51").Value = "  -1.21%  "
